$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("2025-03-13 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-14 Friday", 2)
$null = $d.Content.Find.Execute("43×41=1763", $true, $false, $false, $false, $false, $true, 1, $false, "50×34=1700", 2)
$null = $d.Content.Find.Execute("40×45=1800", $true, $false, $false, $false, $false, $true, 1, $false, "47×93=4371", 2)
$null = $d.Content.Find.Execute("82×91=7462", $true, $false, $false, $false, $false, $true, 1, $false, "52×23=1196", 2)
$null = $d.Content.Find.Execute("74×73=5402", $true, $false, $false, $false, $false, $true, 1, $false, "47×26=1222", 2)
$null = $d.Content.Find.Execute("29×59=1711", $true, $false, $false, $false, $false, $true, 1, $false, "63×82=5166", 2)
$null = $d.Content.Find.Execute("14×61=854", $true, $false, $false, $false, $false, $true, 1, $false, "24×16=384", 2)
$null = $d.Content.Find.Execute("53×98=5194", $true, $false, $false, $false, $false, $true, 1, $false, "94×30=2820", 2)
$null = $d.Content.Find.Execute("31×28=868", $true, $false, $false, $false, $false, $true, 1, $false, "94×12=1128", 2)
$null = $d.Content.Find.Execute("95×98=9310", $true, $false, $false, $false, $false, $true, 1, $false, "74×58=4292", 2)
$null = $d.Content.Find.Execute("72×68=4896", $true, $false, $false, $false, $false, $true, 1, $false, "28×50=1400", 2)
$null = $d.Content.Find.Execute("84×94=7896", $true, $false, $false, $false, $false, $true, 1, $false, "80×37=2960", 2)
$null = $d.Content.Find.Execute("97×64=6208", $true, $false, $false, $false, $false, $true, 1, $false, "63×86=5418", 2)
$null = $d.Content.Find.Execute("62×71=4402", $true, $false, $false, $false, $false, $true, 1, $false, "54×21=1134", 2)
$null = $d.Content.Find.Execute("30×78=2340", $true, $false, $false, $false, $false, $true, 1, $false, "53×40=2120", 2)
$null = $d.Content.Find.Execute("11×60=660", $true, $false, $false, $false, $false, $true, 1, $false, "54×80=4320", 2)
$null = $d.Content.Find.Execute("71×53=3763", $true, $false, $false, $false, $false, $true, 1, $false, "90×88=7920", 2)
$null = $d.Content.Find.Execute("89×37=3293", $true, $false, $false, $false, $false, $true, 1, $false, "14×34=476", 2)
$null = $d.Content.Find.Execute("50×71=3550", $true, $false, $false, $false, $false, $true, 1, $false, "12×86=1032", 2)
$null = $d.Content.Find.Execute("79×73=5767", $true, $false, $false, $false, $false, $true, 1, $false, "18×72=1296", 2)
$null = $d.Content.Find.Execute("85×17=1445", $true, $false, $false, $false, $false, $true, 1, $false, "21×63=1323", 2)
$null = $d.Content.Find.Execute("65×88=5720", $true, $false, $false, $false, $false, $true, 1, $false, "43×45=1935", 2)
$null = $d.Content.Find.Execute("12×94=1128", $true, $false, $false, $false, $false, $true, 1, $false, "47×34=1598", 2)
$null = $d.Content.Find.Execute("80×23=1840", $true, $false, $false, $false, $false, $true, 1, $false, "97×37=3589", 2)
$null = $d.Content.Find.Execute("68×48=3264", $true, $false, $false, $false, $false, $true, 1, $false, "31×17=527", 2)
$null = $d.Content.Find.Execute("60×91=5460", $true, $false, $false, $false, $false, $true, 1, $false, "33×21=693", 2)
